$d = $word.ActiveDocument

# --- Paragraph: 'Programa' (Portuguese syllabus) ---
$find1 = '1. INTRODUÇÃO A CINÉTICATipos de Reações Químicas. Lei de velocidade e seus principais parâmetros. Influência da temperatura sobre a taxa da reação. Ativação das reações químicas Equação de Arrhenius. Energia de ativação. Conversão. Concentração e sua variação numa transformação química. 2. REAÇÕES A VOLUME CONSTANTEReações irreversíveis de ordem um. Reações irreversíveis de ordem dois. Reações irreversíveis de ordem três. Reações irreversíveis de ordem qualquer. 3. REAÇÕES A VOLUME VARIÁVELConceitos. Fração de conversão volumétrica. Reações a volume variável de ordem um e dois. 4. MODELOS IDEAIS DE REATORES QUÍMICOS ISOTÉRMICOS: Equações fundamentais de projeto de reatores. Reator tanque descontínuo (BSTR). Reator tanque de mistura contínuo (CSTR). Reator tubular de fluxo pistonado (PFR). Comparação de desempenho de reatores CSTR e PFR. Reatores CSTR em cascata. Associação mista de reatores em série: CSTR e PFR 5. ANÁLISE DE DADOS CINÉTICOS EM REATORES QUÍMICOS ISOTÉRMICOSBalanço de massa e coleta de dados em reatores ideais isotérmicos: batelada (BSTR), reator tanque de mistura contínuo (CSTR) e Reator tubular (PFR)'
$replace1 = '1. INTRODUÇÃO A CINÉTICA^lTipos de Reações Químicas. Lei de velocidade e seus principais parâmetros. Influência da temperatura sobre a taxa da reação. Ativação das reações químicas Equação de Arrhenius. Energia de ativação. Conversão. Concentração e sua variação numa transformação química. ^l2. REAÇÕES A VOLUME CONSTANTE^lReações irreversíveis de ordem um. Reações irreversíveis de ordem dois. Reações irreversíveis de ordem três. Reações irreversíveis de ordem qualquer. ^l3. REAÇÕES A VOLUME VARIÁVEL^lConceitos. Fração de conversão volumétrica. Reações a volume variável de ordem um e dois. ^l4. MODELOS IDEAIS DE REATORES QUÍMICOS ISOTÉRMICOS: ^lEquações fundamentais de projeto de reatores. Reator tanque descontínuo (BSTR). Reator tanque de mistura contínuo (CSTR). Reator tubular de fluxo pistonado (PFR). Comparação de desempenho de reatores CSTR e PFR. Reatores CSTR em cascata. Associação mista de reatores em série: CSTR e PFR ^l5. ANÁLISE DE DADOS CINÉTICOS EM REATORES QUÍMICOS ISOTÉRMICOS^lBalanço de massa e coleta de dados em reatores ideais isotérmicos: batelada (BSTR), reator tanque de mistura contínuo (CSTR) e Reator tubular (PFR)'
$found1 = $d.Content.Find.Execute($find1, $false, $false, $false, $false, $false, $true, 1, $false, $replace1, 2)
Write-Host "Programa (PT): $found1"

# --- Paragraph: 'Programa' (English syllabus, italic run) ---
$find2 = '1. Introduction to KineticsTypes of Chemical Reactions.Rate law and its main parameters.Influence of temperature on reaction rate.Activation of chemical reactions.Arrhenius equation.Activation energy.Conversion.Concentration and its variation in a chemical transformation. 2. Reactions at Constant VolumeIrreversible reactions of first order.Irreversible reactions of second order.Irreversible reactions of third order.Irreversible reactions of any order. 3. Reactions at Variable VolumeConcepts.Volumetric conversion fraction.Reactions at variable volume of first and second order. 4. Ideal Models of Isothermal Chemical Reactors:Fundamental equations for reactor design.Batch reactor (BSTR).Continuous stirred-tank reactor (CSTR).Plug-flow reactor (PFR).Performance comparison of CSTR and PFR.Cascade CSTR reactors.Mixed association of reactors in series: CSTR and PFR. 5. Analysis of Kinetic Data in Isothermal Chemical ReactorsMass balance and data collection in ideal isothermal reactors:Batch reactor (BSTR).Continuous stirred-tank reactor (CSTR).Plug-flow reactor (PFR).'
$replace2 = '1. Introduction to Kinetics^lTypes of Chemical Reactions.^lRate law and its main parameters.^lInfluence of temperature on reaction rate.^lActivation of chemical reactions.^lArrhenius equation.^lActivation energy.^lConversion.^lConcentration and its variation in a chemical transformation. ^l2. Reactions at Constant Volume^lIrreversible reactions of first order.^lIrreversible reactions of second order.^lIrreversible reactions of third order.^lIrreversible reactions of any order. ^l3. Reactions at Variable Volume^lConcepts.^lVolumetric conversion fraction.^lReactions at variable volume of first and second order. ^l4. Ideal Models of Isothermal Chemical Reactors:^lFundamental equations for reactor design.^lBatch reactor (BSTR).^lContinuous stirred-tank reactor (CSTR).^lPlug-flow reactor (PFR).^lPerformance comparison of CSTR and PFR.^lCascade CSTR reactors.^lMixed association of reactors in series: CSTR and PFR. ^l5. Analysis of Kinetic Data in Isothermal Chemical Reactors^lMass balance and data collection in ideal isothermal reactors:^lBatch reactor (BSTR).^lContinuous stirred-tank reactor (CSTR).^lPlug-flow reactor (PFR).'
$found2 = $d.Content.Find.Execute($find2, $false, $false, $false, $false, $false, $true, 1, $false, $replace2, 2)
Write-Host "Programa (EN): $found2"

# --- Paragraph: 'Criterio' sentence inside Avaliacao ---
$find3 = 'Média da Primeira Avaliação (N) = 50% P1 + 50% P2.Obs: fica a critério de cada docente a inserção de trabalhos no decorrer do curso, bem como a alteração do peso de cada prova em decorrência dos mesmos.'
$replace3 = 'Média da Primeira Avaliação (N) = 50% P1 + 50% P2.^lObs: fica a critério de cada docente a inserção de trabalhos no decorrer do curso, bem como a alteração do peso de cada prova em decorrência dos mesmos.'
$found3 = $d.Content.Find.Execute($find3, $false, $false, $false, $false, $false, $true, 1, $false, $replace3, 2)
Write-Host "Criterio: $found3"

# --- Paragraph: 'Bibliografia' ---
$find4 = '1- FOGLER, H.S. Elementos de engenharia das reações químicas. 3.ed. Rio de Janeiro: LTC Editora, 2009.2- LEVENSPIEL, O. Engenharia Das Reações Químicas, E ed (Blucher, São Paulo, 2000)3- VAN SANTEN, R.A.; Niemantsverdriet, J.W. Chemical kinetics and catalysis. New York: Plenum Press, 1995.4- Missen, R.W.; Mims, C.A.; Saville, B.A. Introduction to chemical reaction engineering and kinetics. New York: J. Wiley, 1999.5- Rothenberg, G. Catalysis: concepts and green applications. Weinheim: Wiley-VCH, 2008 Chichester.6- DENISOV, E.T.; Sarkisov, O.M.; Likhtenshtein, G.I. Chemical kinetics: fundamentals and new developments. Amsterdam: Elsevier, 2003.7- Hagen, J. Industrial catalysis: a practical approach. Weinheim: Wiley-VCH, 2006.8- Salmi, T.O.; Mikkola, J.; Warna, J.P. Chemical reaction engineering and reactor technology. Boca Raton: CRC Press/Taylor & Francis, 2011.9- Mortimer, M.; Taylor, P.G. Chemical kinetics and mechanism. Cambridge: Royal Society of Chemistry, 2002.10- FROMENT, G.F.; BISCHOFF, K.B. Chemical reactor analysis and design. 2nd. Ed. New York: John Wiley & Sons, 1990.11- HILL, C.G. An Introduction to chemical engineering kinetics and reactor design. New York: John Wiley&Sons, 1977.12- SMITH, J.M. Chemical engineering kinetics. 3rd. ed New York: McGraw-Hill,1981.13- DENBIGH, K.; TURNER, R. Introduction to chemical Reaction Design. Cambridge: Cambridge University Press, 1970.14 - AGUIAR, L. G. Problemas de cinética e reatores químicos. Curitiba: Appris Editora, 2023.'
$replace4 = '1- FOGLER, H.S. Elementos de engenharia das reações químicas. 3.ed. Rio de Janeiro: LTC Editora, 2009.^l^l2- LEVENSPIEL, O. Engenharia Das Reações Químicas, E ed (Blucher, São Paulo, 2000)^l3- VAN SANTEN, R.A.; Niemantsverdriet, J.W. Chemical kinetics and catalysis. New York: Plenum Press, 1995.^l4- Missen, R.W.; Mims, C.A.; Saville, B.A. Introduction to chemical reaction engineering and kinetics. New York: J. Wiley, 1999.^l5- Rothenberg, G. Catalysis: concepts and green applications. Weinheim: Wiley-VCH, 2008 Chichester.^l6- DENISOV, E.T.; Sarkisov, O.M.; Likhtenshtein, G.I. Chemical kinetics: fundamentals and new developments. Amsterdam: Elsevier, 2003.^l7- Hagen, J. Industrial catalysis: a practical approach. Weinheim: Wiley-VCH, 2006.^l8- Salmi, T.O.; Mikkola, J.; Warna, J.P. Chemical reaction engineering and reactor technology. Boca Raton: CRC Press/Taylor & Francis, 2011.^l9- Mortimer, M.; Taylor, P.G. Chemical kinetics and mechanism. Cambridge: Royal Society of Chemistry, 2002.^l10- FROMENT, G.F.; BISCHOFF, K.B. Chemical reactor analysis and design. 2nd. Ed. New York: John Wiley & Sons, 1990.^l11- HILL, C.G. An Introduction to chemical engineering kinetics and reactor design. New York: John Wiley&Sons, 1977.^l12- SMITH, J.M. Chemical engineering kinetics. 3rd. ed New York: McGraw-Hill,1981.^l13- DENBIGH, K.; TURNER, R. Introduction to chemical Reaction Design. Cambridge: Cambridge University Press, 1970.^l14 - AGUIAR, L. G. Problemas de cinética e reatores químicos. Curitiba: Appris Editora, 2023.'
$found4 = $d.Content.Find.Execute($find4, $false, $false, $false, $false, $false, $true, 1, $false, $replace4, 2)
Write-Host "Bibliografia: $found4"
